$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells per diff (rows 2,3,4,6,7,8,10,11,13) ---
$ws.Range("F2").Value = 4.3
$ws.Range("G2").Value = 6.2
$ws.Range("H2").Value = 1.71
$ws.Range("I2").Value = 1.95
$ws.Range("J2").Value = 3.55
$ws.Range("K2").Value = 4.9
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 3.75
$ws.Range("O2").Value = 1.26
$ws.Range("P2").Value = 2.1
$ws.Range("R2").Value = 1.43
$ws.Range("S2").Value = 2.76
$ws.Range("V2").Value = 2.04
$ws.Range("W2").Value = 1.19
$ws.Range("T3").Value = 2.32
$ws.Range("U3").Value = 1.65
$ws.Range("G4").Value = 1.8
$ws.Range("N4").Value = 3.8
$ws.Range("O4").Value = 1.28
$ws.Range("P4").Value = 1.97
$ws.Range("R4").Value = 1.37
$ws.Range("S4").Value = 3.15
$ws.Range("T4").Value = 1.83
$ws.Range("U4").Value = 1.98
$ws.Range("W4").Value = 2.24
$ws.Range("X4").Value = 18.5
$ws.Range("Y4").Value = 23
$ws.Range("Z4").Value = 55
$ws.Range("AB4").Value = 10
$ws.Range("AC4").Value = 11
$ws.Range("AD4").Value = 26
$ws.Range("AF4").Value = 12
$ws.Range("AG4").Value = 11.5
$ws.Range("AH4").Value = 24
$ws.Range("AJ4").Value = 21
$ws.Range("AK4").Value = 21
$ws.Range("AL4").Value = 42
$ws.Range("AN4").Value = 12
$ws.Range("L6").Value = 1.19
$ws.Range("P6").Value = 3
$ws.Range("Q6").Value = 1.39
$ws.Range("R6").Value = 1.81
$ws.Range("U6").Value = 1.85
$ws.Range("X6").Value = 46
$ws.Range("Z6").Value = 180
$ws.Range("AB6").Value = 14.5
$ws.Range("AC6").Value = 21
$ws.Range("AD6").Value = 65
$ws.Range("AE6").Value = 250
$ws.Range("AF6").Value = 11
$ws.Range("AG6").Value = 14
$ws.Range("AH6").Value = 40
$ws.Range("AI6").Value = 180
$ws.Range("AJ6").Value = 11.5
$ws.Range("AK6").Value = 16
$ws.Range("AL6").Value = 42
$ws.Range("AM6").Value = 170
$ws.Range("AN6").Value = 3.75
$ws.Range("U7").Value = 1.78
$ws.Range("Z7").Value = 10
$ws.Range("P8").Value = 1.38
$ws.Range("U8").Value = 1.61
$ws.Range("V8").Value = 1.39
$ws.Range("W8").Value = 1.51
$ws.Range("Y8").Value = 8.199999999999999
$ws.Range("Z8").Value = 22
$ws.Range("AA8").Value = 100
$ws.Range("AB8").Value = 7.2
$ws.Range("AD8").Value = 17
$ws.Range("AF8").Value = 19.5
$ws.Range("AG8").Value = 18.5
$ws.Range("AI8").Value = 140
$ws.Range("AJ8").Value = 55
$ws.Range("AK8").Value = 55
$ws.Range("AN8").Value = 85
$ws.Range("AO8").Value = 140
$ws.Range("L10").Value = 1.47
$ws.Range("F11").Value = 2.5
$ws.Range("G11").Value = 2.82
$ws.Range("H11").Value = 3.2
$ws.Range("I11").Value = 3.7
$ws.Range("J11").Value = 2.84
$ws.Range("N11").Value = 2.48
$ws.Range("V11").Value = 1.37
$ws.Range("W11").Value = 1.54
$ws.Range("X11").Value = 9.800000000000001
$ws.Range("Y11").Value = 10
$ws.Range("AA11").Value = 90
$ws.Range("AB11").Value = 9.800000000000001
$ws.Range("AE11").Value = 55
$ws.Range("AF11").Value = 21
$ws.Range("AI11").Value = 100
$ws.Range("AJ11").Value = 55
$ws.Range("AK11").Value = 46
$ws.Range("AN11").Value = 55
$ws.Range("AO11").Value = 95
$ws.Range("G13").Value = 2.52
$ws.Range("I13").Value = 4.4
$ws.Range("J13").Value = 3.05
$ws.Range("S13").Value = 4.4

# --- Add new rows 15 and 16 (new fixtures appended below existing data) ---
# Pre-format the Date column (B) for the new rows as Text so the
# "YYYY-MM-DD" strings are stored literally instead of being auto-converted
# into date serial numbers (matches the Date formatting used on import).
$ws.Range("B15:B16").NumberFormat = "@"

# Row 15
$ws.Range("A15").Value = "Brazilian Serie A"
$ws.Range("B15").Value = "2026-02-05"
$ws.Range("C15").Value = "21:30:00"
$ws.Range("D15").Value = "Cruzeiro MG"
$ws.Range("E15").Value = "Coritiba"
$ws.Range("F15").Value = 1.45
$ws.Range("G15").Value = 1.5
$ws.Range("H15").Value = 8.800000000000001
$ws.Range("I15").Value = 10.5
$ws.Range("J15").Value = 4.3
$ws.Range("K15").Value = 4.9
$ws.Range("L15").Value = 1.4
$ws.Range("M15").Value = 1.07
$ws.Range("N15").Value = 3.35
$ws.Range("O15").Value = 1.36
$ws.Range("P15").Value = 1.82
$ws.Range("Q15").Value = 2.06
$ws.Range("R15").Value = 1.31
$ws.Range("S15").Value = 3.75
$ws.Range("T15").Value = 2.26
$ws.Range("U15").Value = 1.68
$ws.Range("V15").Value = 1.11
$ws.Range("W15").Value = 2.96
$ws.Range("X15").Value = 970
$ws.Range("Y15").Value = 970
$ws.Range("Z15").Value = 110
$ws.Range("AA15").Value = 490
$ws.Range("AB15").Value = 8
$ws.Range("AC15").Value = 11
$ws.Range("AD15").Value = 970
$ws.Range("AE15").Value = 230
$ws.Range("AF15").Value = 7.8
$ws.Range("AG15").Value = 970
$ws.Range("AH15").Value = 970
$ws.Range("AI15").Value = 200
$ws.Range("AJ15").Value = 13
$ws.Range("AK15").Value = 970
$ws.Range("AL15").Value = 50
$ws.Range("AM15").Value = 280
$ws.Range("AN15").Value = 9.4
$ws.Range("AO15").Value = 410
# Row 16
$ws.Range("A16").Value = "Colombian Primera A"
$ws.Range("B16").Value = "2026-02-05"
$ws.Range("C16").Value = "22:00:00"
$ws.Range("D16").Value = "Millonarios"
$ws.Range("E16").Value = "Deportivo Pereira"
$ws.Range("F16").Value = 1.47
$ws.Range("G16").Value = 1.53
$ws.Range("H16").Value = 9.4
$ws.Range("I16").Value = 12
$ws.Range("J16").Value = 3.85
$ws.Range("K16").Value = 4.6
$ws.Range("L16").Value = 1.01
$ws.Range("M16").Value = 1.1
$ws.Range("N16").Value = 2.8
$ws.Range("O16").Value = 1.49
$ws.Range("P16").Value = 1.6
$ws.Range("Q16").Value = 2.46
$ws.Range("R16").Value = 1.2
$ws.Range("S16").Value = 4.9
$ws.Range("T16").Value = 2.52
$ws.Range("U16").Value = 1.54
$ws.Range("V16").Value = 1.09
$ws.Range("W16").Value = 2.88
$ws.Range("X16").Value = 11
$ws.Range("Y16").Value = 24
$ws.Range("Z16").Value = 1000
$ws.Range("AA16").Value = 1000
$ws.Range("AB16").Value = 5.7
$ws.Range("AC16").Value = 10.5
$ws.Range("AD16").Value = 970
$ws.Range("AE16").Value = 340
$ws.Range("AF16").Value = 7.4
$ws.Range("AG16").Value = 11.5
$ws.Range("AH16").Value = 970
$ws.Range("AI16").Value = 320
$ws.Range("AJ16").Value = 14
$ws.Range("AK16").Value = 23
$ws.Range("AL16").Value = 75
$ws.Range("AM16").Value = 470
$ws.Range("AN16").Value = 13.5
$ws.Range("AO16").Value = 1000
